# Add Test Data for UK Market:
# Duplicate the "Poland" sheet (the existing last country-market template
# sheet) to create a new "UK" sheet at the end of the workbook, then fill
# in its Market / Jira-ticket cells and make it the active sheet/tab.

$wb = $excel.ActiveWorkbook

$poland = $wb.Worksheets.Item("Poland")

# Copy Poland to a position right after itself (i.e. appended at the end
# of the workbook, after the last existing sheet).
[void]$poland.Copy($null, $poland)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "UK"

# B4 holds the Jira/NGC ticket reference, B2 holds the "<Country> Market"
# label - set B4 first so the new shared strings are appended in the same
# order the original authoring tool produced them (ticket, then market).
$newSheet.Range("B4").Value = "NGC-2741/T3344"
$newSheet.Range("B2").Value = "UK Market"

# Make the new sheet the active tab, with B4 selected (matching the state
# left behind right after typing the ticket reference into B4).
[void]$newSheet.Activate()
[void]$newSheet.Range("B4").Select()
